$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7
$ws.Range("B7").Value = 0.06372808675650223
$ws.Range("C7").Value = 0.9343998086003413
$ws.Range("D7").Value = 2.775446151470913
$ws.Range("E7").Value = 1.665967031927977
$ws.Range("F7").Value = 1.688362315275229
$ws.Range("G7").Value = 36

# Row 8
$ws.Range("B8").Value = 0.0399851282911067
$ws.Range("C8").Value = 1.050477138034481
$ws.Range("D8").Value = 3.029802060013869
$ws.Range("E8").Value = 1.740632660848885
$ws.Range("F8").Value = 1.765578673122826
$ws.Range("G8").Value = 35

# Row 9
$ws.Range("B9").Value = -0.04534382660436193
$ws.Range("C9").Value = 1.253794741722951
$ws.Range("D9").Value = 4.508050214875826
$ws.Range("E9").Value = 2.123216949554573
$ws.Range("F9").Value = 2.177877806311915
$ws.Range("G9").Value = 20

# Row 10
$ws.Range("B10").Value = -0.4013436200723844
$ws.Range("C10").Value = 0.9917379900131201
$ws.Range("D10").Value = 2.080463190564922
$ws.Range("E10").Value = 1.442381083682437
$ws.Range("F10").Value = 1.441990070455223
$ws.Range("G10").Value = 13

# Row 11
$ws.Range("B11").Value = 0.2679351619938605
$ws.Range("C11").Value = 0.5450213326915778
$ws.Range("D11").Value = 0.3698384655848825
$ws.Range("E11").Value = 0.6081434580630483
$ws.Range("F11").Value = 0.6103781763712214
$ws.Range("G11").Value = 5
